$wb = $excel.ActiveWorkbook

# Sheet "2025": update base value in A2 and B2
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("A2").Value = 324500
$ws2025.Range("B2").Value = 11

# Sheet "2030": A2 becomes a formula referencing '2025'!A2, B2 updated
$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Formula = "='2025'!A2*(1-0.15*0.2)"
$ws2030.Range("B2").Value = 11

# Sheet "2035": A2 becomes a formula referencing '2025'!A2, B2 updated
$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Formula = "='2025'!A2*(1-0.15*0.4)"
$ws2035.Range("B2").Value = 11

# Sheet "2040": A2 becomes a formula referencing '2025'!A2, B2 updated
$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("A2").Formula = "='2025'!A2*(1-0.15*0.6)"
$ws2040.Range("B2").Value = 11

# Sheet "2045": A2 becomes a formula referencing '2025'!A2, B2 updated
$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Formula = "='2025'!A2*(1-0.15*0.8)"
$ws2045.Range("B2").Value = 11

# Sheet "2050": A2 becomes a formula referencing '2025'!A2, B2 updated
$ws2050 = $wb.Worksheets.Item("2050")
$ws2050.Range("A2").Formula = "='2025'!A2*(1-0.15*1)"
$ws2050.Range("B2").Value = 11

$wb.Save()
